# Fixed naive component forecaster bug - Presentation state 11.02.
# Applies corrected YoY forecast values (re-run with fixed naive-component
# forecaster) and removes a stray duplicate value in C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").ClearContents()
$ws.Range("E2").Value = 5.080273296954396
$ws.Range("C3").Value = -3.942037578692492
$ws.Range("E3").Value = -1.648748515828491
$ws.Range("C4").Value = -2.839753013810498
$ws.Range("E4").Value = -1.632723506456923
$ws.Range("C5").Value = 4.960109259035428
$ws.Range("C7").Value = -2.700325749999499
$ws.Range("E7").Value = -0.3858735870725494
$ws.Range("C8").Value = 5.469647210234974
$ws.Range("E8").Value = 3.061326532789521
$ws.Range("C10").Value = 3.458696398997096
$ws.Range("E10").Value = 2.610227683091315
$ws.Range("C11").Value = 2.77241330895972
$ws.Range("C12").Value = 3.14581984265847
$ws.Range("E13").Value = 4.124307769579483
$ws.Range("E14").Value = 4.888255652935958
$ws.Range("E15").Value = 2.95288809451808
$ws.Range("C18").Value = -0.5744163079740128
$ws.Range("E18").Value = -0.6956477387308979
$ws.Range("C19").Value = -0.1892239049850142
